$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new exam row (2021 - Vår) below the existing table (row 27)
$ws.Range("A27").Value = "2021 - Vår"
$ws.Range("B27").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-21-v.pdf)"
$ws.Range("C27").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-21-v-fasit.pdf)"

# Leave the selection where the user ended up after adding the row
$ws.Range("C28").Select()
